# Auto-generated edit script applying numeric updates to the Leve profit sheets
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ /
#  LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns H..N), matching the scheduled
# market-data runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3479
$ws.Range("I64").Value = 2895
$ws.Range("J64").Value = 3625
$ws.Range("K64").Value = 2895
$ws.Range("L64").Value = 3625
$ws.Range("M64").Value = -2647
$ws.Range("N64").Value = -4121
$ws.Range("H67").Value = 3479
$ws.Range("I67").Value = 2895
$ws.Range("J67").Value = 3625
$ws.Range("K67").Value = 2895
$ws.Range("L67").Value = 3625
$ws.Range("M67").Value = -2037
$ws.Range("N67").Value = -5341
$ws.Range("H92").Value = 793.4
$ws.Range("I92").Value = 793.4
$ws.Range("K92").Value = 793.4
$ws.Range("M92").Value = 454.6
$ws.Range("H106").Value = 3368.7407
$ws.Range("I106").Value = 2957
$ws.Range("K106").Value = 2957
$ws.Range("M106").Value = -2326
$ws.Range("H113").Value = 58827268
$ws.Range("J113").Value = 5180.4
$ws.Range("L113").Value = 5180.4
$ws.Range("N113").Value = -11688.4
$ws.Range("H129").Value = 257331.31
$ws.Range("J129").Value = 278753.44
$ws.Range("L129").Value = 836260.3200000001
$ws.Range("N129").Value = -846260.3200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10357.486
$ws.Range("I32").Value = 7912.459
$ws.Range("J32").Value = 23916.273
$ws.Range("K32").Value = 7912.459
$ws.Range("L32").Value = 23916.273
$ws.Range("M32").Value = -7625.459
$ws.Range("N32").Value = -24490.273
$ws.Range("H74").Value = 31251628
$ws.Range("I74").Value = 43478976
$ws.Range("J74").Value = 3955.5557
$ws.Range("K74").Value = 43478976
$ws.Range("L74").Value = 3955.5557
$ws.Range("M74").Value = -43478102
$ws.Range("N74").Value = -5703.5557
$ws.Range("H77").Value = 31251628
$ws.Range("I77").Value = 43478976
$ws.Range("J77").Value = 3955.5557
$ws.Range("K77").Value = 217394880
$ws.Range("L77").Value = 19777.7785
$ws.Range("M77").Value = -217390512
$ws.Range("N77").Value = -28513.7785
$ws.Range("H132").Value = 9627141
$ws.Range("I132").Value = 11629798
$ws.Range("J132").Value = 58893.223
$ws.Range("K132").Value = 34889394
$ws.Range("L132").Value = 176679.669
$ws.Range("M132").Value = -34886864
$ws.Range("N132").Value = -181739.669

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1170.9231
$ws.Range("I80").Value = 1002.4167
$ws.Range("J80").Value = 1315.3572
$ws.Range("K80").Value = 1002.4167
$ws.Range("L80").Value = 1315.3572
$ws.Range("M80").Value = -4.416699999999992
$ws.Range("N80").Value = -3311.3572
$ws.Range("H83").Value = 1170.9231
$ws.Range("I83").Value = 1002.4167
$ws.Range("J83").Value = 1315.3572
$ws.Range("K83").Value = 5012.0835
$ws.Range("L83").Value = 6576.786
$ws.Range("M83").Value = -20.08349999999973
$ws.Range("N83").Value = -16560.786
$ws.Range("H94").Value = 1107.2916
$ws.Range("I94").Value = 611.9286
$ws.Range("K94").Value = 611.9286
$ws.Range("M94").Value = -160.9286
$ws.Range("H99").Value = 1387.7778
$ws.Range("I99").Value = 1423.75
$ws.Range("K99").Value = 1423.75
$ws.Range("M99").Value = 74.25
$ws.Range("H105").Value = 1615220
$ws.Range("I105").Value = 1777.5
$ws.Range("J105").Value = 2176417.2
$ws.Range("K105").Value = 1777.5
$ws.Range("L105").Value = 2176417.2
$ws.Range("M105").Value = -30.5
$ws.Range("N105").Value = -2179911.2
$ws.Range("H126").Value = 59999
$ws.Range("J126").Value = 59999
$ws.Range("L126").Value = 59999
$ws.Range("N126").Value = -69879

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12707.863
$ws.Range("I58").Value = 1482.2222
$ws.Range("J58").Value = 20479.46
$ws.Range("K58").Value = 1482.2222
$ws.Range("L58").Value = 20479.46
$ws.Range("M58").Value = -1279.2222
$ws.Range("N58").Value = -20885.46
$ws.Range("H81").Value = 38000
$ws.Range("J81").Value = 38000
$ws.Range("L81").Value = 38000
$ws.Range("N81").Value = -39996
$ws.Range("H84").Value = 38000
$ws.Range("J84").Value = 38000
$ws.Range("L84").Value = 114000
$ws.Range("N84").Value = -123984
$ws.Range("H86").Value = 9697.4
$ws.Range("I86").Value = 1923.375
$ws.Range("J86").Value = 18582
$ws.Range("K86").Value = 1923.375
$ws.Range("L86").Value = 18582
$ws.Range("M86").Value = -800.375
$ws.Range("N86").Value = -20828
$ws.Range("H89").Value = 9697.4
$ws.Range("I89").Value = 1923.375
$ws.Range("J89").Value = 18582
$ws.Range("K89").Value = 9616.875
$ws.Range("L89").Value = 92910
$ws.Range("M89").Value = -4000.875
$ws.Range("N89").Value = -104142
$ws.Range("H99").Value = 3515.2144
$ws.Range("I99").Value = 2639.818
$ws.Range("J99").Value = 6725
$ws.Range("K99").Value = 2639.818
$ws.Range("L99").Value = 6725
$ws.Range("M99").Value = -1141.818
$ws.Range("N99").Value = -9721
$ws.Range("H122").Value = 1240.6072
$ws.Range("I122").Value = 1030.3334
$ws.Range("K122").Value = 3091.0002
$ws.Range("M122").Value = -641.0001999999999
$ws.Range("H126").Value = 3515.2144
$ws.Range("I126").Value = 2639.818
$ws.Range("J126").Value = 6725
$ws.Range("K126").Value = 7919.454000000001
$ws.Range("L126").Value = 20175
$ws.Range("M126").Value = -5449.454000000001
$ws.Range("N126").Value = -25115
$ws.Range("H136").Value = 12707.863
$ws.Range("I136").Value = 1482.2222
$ws.Range("J136").Value = 20479.46
$ws.Range("K136").Value = 4446.6666
$ws.Range("L136").Value = 61438.38
$ws.Range("M136").Value = -1896.6666
$ws.Range("N136").Value = -66538.38

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 797.9091
$ws.Range("J113").Value = 883.1539
$ws.Range("L113").Value = 2649.4617
$ws.Range("N113").Value = -6989.4617
$ws.Range("H121").Value = 2372.6
$ws.Range("J121").Value = 3777.6667
$ws.Range("L121").Value = 11333.0001
$ws.Range("N121").Value = -13953.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3680897.2
$ws.Range("I70").Value = 4714.143
$ws.Range("K70").Value = 4714.143
$ws.Range("M70").Value = -4444.143
$ws.Range("H73").Value = 3680897.2
$ws.Range("I73").Value = 4714.143
$ws.Range("K73").Value = 4714.143
$ws.Range("M73").Value = -3778.143
$ws.Range("H102").Value = 4934.8
$ws.Range("I102").Value = 4779.25
$ws.Range("K102").Value = 4779.25
$ws.Range("M102").Value = -3157.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3067.125
$ws.Range("I7").Value = 3027.4
$ws.Range("J7").Value = 3133.3333
$ws.Range("K7").Value = 3027.4
$ws.Range("L7").Value = 3133.3333
$ws.Range("M7").Value = -2915.4
$ws.Range("N7").Value = -3357.3333
$ws.Range("H40").Value = 97936.914
$ws.Range("I40").Value = 115924.3
$ws.Range("K40").Value = 115924.3
$ws.Range("M40").Value = -115788.3
$ws.Range("H122").Value = 1637044.4
$ws.Range("I122").Value = 1963153.2
$ws.Range("K122").Value = 5889459.6
$ws.Range("M122").Value = -5887009.6
$ws.Range("H126").Value = 3067.125
$ws.Range("I126").Value = 3027.4
$ws.Range("J126").Value = 3133.3333
$ws.Range("K126").Value = 9082.200000000001
$ws.Range("L126").Value = 9399.999899999999
$ws.Range("M126").Value = -6612.200000000001
$ws.Range("N126").Value = -14339.9999
$ws.Range("H136").Value = 2740.55
$ws.Range("I136").Value = 2740.55
$ws.Range("K136").Value = 8221.650000000001
$ws.Range("M136").Value = -5671.650000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 4660.4
$ws.Range("I4").Value = 3434
$ws.Range("J4").Value = 6500
$ws.Range("K4").Value = 3434
$ws.Range("L4").Value = 6500
$ws.Range("M4").Value = -3321
$ws.Range("N4").Value = -6726
$ws.Range("H126").Value = 2118.6
$ws.Range("I126").Value = 1754.7858
$ws.Range("K126").Value = 5264.357400000001
$ws.Range("M126").Value = -2794.357400000001
$ws.Range("H136").Value = 33336916
$ws.Range("I136").Value = 52633550
$ws.Range("J136").Value = 6364.091
$ws.Range("K136").Value = 157900650
$ws.Range("L136").Value = 19092.273
$ws.Range("M136").Value = -157898100
$ws.Range("N136").Value = -24192.273

Write-Host "Applied Typhon_Profits market-data updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets"
